$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 119.
# This shifts rows 119..250 down to 120..251 (so the former last row, 250,
# becomes the new row 251), and leaves a blank row 119 ready to be filled
# in with the new record below.
$ws.Rows("119:119").Insert()

# Populate the new row 119 with the new weekly record.
$row = 119
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = "2021-12-09"
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112023
$ws.Cells.Item($row, 7).Value = "Brócoli"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 500
$ws.Cells.Item($row, 11).Value = 1000
$ws.Cells.Item($row, 12).Value = 1200
$ws.Cells.Item($row, 13).Value = 1100
$ws.Cells.Item($row, 14).Value = "`$/unidad"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1100
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
